# Insert a new data row at row 646 (shifting existing rows 646:687 down to
# 647:688) and populate it with the new daily entry
# (2026/01/14, 水, 7, 201). This mirrors the upstream diff exactly: it is a
# pure row-insert, every row below shifts down by one, and a single brand
# new row of data appears at the bottom (old row 687 -> new row 688).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 646, pushing 646:687 -> 647:688.
$ws.Rows(646).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/14"), not a
# real Excel date. Leading apostrophe forces text entry so Excel does not
# silently convert it to a date serial number; resetting the Style to
# "Normal" afterwards drops the quote-prefix style so the cell ends up with
# no explicit style, matching its sibling rows.
$ws.Range("A646").Value = "'2026/01/14"
$ws.Range("A646").Style = "Normal"

$ws.Range("B646").Value = "水"
$ws.Range("C646").Value = 7
$ws.Range("D646").Value = 201
